# Generate Report for Handback
# Updates timestamps / status values recorded by the handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the two
# rows that shared the old timestamp. The de-de sheet's "Correspond
# Handoff Datetime" column (H) happened to hold the exact same timestamp
# text for these rows, so it is refreshed together with it.
$wsOverview.Range("G2").Value = "2016-08-16 20:13:52"
$wsOverview.Range("G3").Value = "2016-08-16 20:13:52"
$wsDeDe.Range("H2").Value = "2016-08-16 20:13:52"
$wsDeDe.Range("H3").Value = "2016-08-16 20:13:52"

# zh-cn sheet: Priority column (E) status changed from "ht" to "mt".
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback
# DateTime (K) refreshed to the new handback run's timestamps.
$wsZhCn.Range("H2").Value = "2016-08-16 20:13:45"
$wsZhCn.Range("H3").Value = "2016-08-16 20:13:45"
$wsZhCn.Range("K2").Value = "2016-08-16 20:14:17"
$wsZhCn.Range("K3").Value = "2016-08-16 20:14:17"

# de-de sheet: Priority column (E) status changed from "ht" to "mt".
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# de-de sheet: Correspond Handback DateTime (K) refreshed.
$wsDeDe.Range("K2").Value = "2016-08-16 20:14:24"
$wsDeDe.Range("K3").Value = "2016-08-16 20:14:24"
